# Disposable reverse proxy architecture - update demo deck
# Applies:
#  - date placeholder text bump (03/01/2024 -> 04/01/2024) on master + all layouts
#  - "Sailing servers (...)" flowchart box: taller + reworded detail text
#  - "Temp rules without dns/route53 rule" box: taller + reworded text
#  - auto-routed connector between those two boxes: updated geometry
#  - "ALB are separate listeners" footnote: pluralised to "ALBs"

function EmuToPt($emu) {
    # PowerPoint COM geometry is expressed in points (1 pt = 12700 EMU).
    # Nudge by half an EMU so the float -> EMU round-trip on save lands on
    # the exact target EMU value instead of truncating one unit short.
    return ($emu / 12700.0) + (0.5 / 12700.0)
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Date placeholder fields: 03/01/2024 -> 04/01/2024
#    (slide master + every slide layout)
# ---------------------------------------------------------------------
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        if ($shp.TextFrame.TextRange.Text -eq "03/01/2024") {
            $shp.TextFrame.TextRange.Text = "04/01/2024"
        }
    }
}

$layouts = $master.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    $layout = $layouts.Item($L)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq "03/01/2024") {
                $shp.TextFrame.TextRange.Text = "04/01/2024"
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2. Slide 1 shape edits
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)

# --- "Flowchart: Multidocument 54" (Sailing servers box) ---
$sailingBox = $s.Shapes.Item(18)
$sailingBox.Height = EmuToPt 1259870
$tr = $sailingBox.TextFrame.TextRange
$tr.Characters(17, 89).Text = "(be that on a server, an independent instance, perhaps some replicas the master is independent but the replicas are shared, etc.)"

# --- "TextBox 58" (Temp rules without dns/route53 rule) ---
$tempRulesBox = $s.Shapes.Item(20)
$tempRulesBox.Height = EmuToPt 600164
$tr2 = $tempRulesBox.TextFrame.TextRange
$tr2.Characters(23, 13).Text = "/route53 rule, which point to sailing servers"
$tr2.Characters(1, 19).Text = "Temporary rules without "

# --- "Straight Arrow Connector 64" (auto-routed between the two boxes above) ---
$connector = $s.Shapes.Item(22)
$connector.Left = EmuToPt 10130118
$connector.Top = EmuToPt 2907887
$connector.Width = EmuToPt 328779
$connector.Height = EmuToPt 2199567

# --- "TextBox 105" (footnote: ALB -> ALBs) ---
$footnoteBox = $s.Shapes.Item(38)
$tr3 = $footnoteBox.TextFrame.TextRange
$tr3.Characters(126, 27).Text = " ALBs are separate listeners"
